# "Admin Request status update capability"
#
# The sessions template's footer note ("Please do not change the column
# order.") used to live down in H2, leaving an orphan row under the
# single-row header. Promote it up into row 1 (H1), alongside the rest of
# the header cells, and drop the now-empty second row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy H2 -> H1 using Range.Copy(Destination) so the cell's style (fill,
# font, number format, etc.) comes along with the value instead of just
# the raw content.
$ws.Range("H2").Copy($ws.Range("H1"))

# H2 (and the row it was alone in) is no longer needed.
$ws.Range("H2").Clear()

# Match the author's resulting selection/view: H1 is the active cell.
$ws.Range("H1").Select()
